# Auto-generated Excel COM-interop script
# Re-applies a fresh nonlinear-regression pass: updates the fitted
# constants/std-deviations, the adjusted R^2, the per-point equilibrium
# concentrations, and the per-titration-point EMF absolute/relative errors.
$wb = $excel.ActiveWorkbook

# --- constants_evaluated: Constant/St.Deviation text cells for component "HL" ---
$ws = $wb.Worksheets.Item("constants_evaluated")
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "4.97236251831055"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "0.0185475079135167"

# --- adj_r_squared ---
$ws = $wb.Worksheets.Item("adj_r_squared")
$ws.Range("A2").Value = [double]"0.981565247999872"

# --- equilibrium_concentrations: A2:F11 ---
$ws = $wb.Worksheets.Item("equilibrium_concentrations")
$ws.Range("A2").Value = [double]"0.000320489792950694"
$ws.Range("B2").Value = [double]"0.000327293669620107"
$ws.Range("C2").Value = [double]"2.80236175548005e-16"
$ws.Range("D2").Value = [double]"0.00984270176871473"
$ws.Range("E2").Value = [double]"4.56166515283288e-09"
$ws.Range("F2").Value = [double]"6.81300017301839e-06"
$ws.Range("A3").Value = [double]"0.000144815394409726"
$ws.Range("B3").Value = [double]"0.000694372558467839"
$ws.Range("C3").Value = [double]"5.00301728233769e-14"
$ws.Range("D3").Value = [double]"0.00943560602358625"
$ws.Range("E3").Value = [double]"2.14179459284501e-08"
$ws.Range("F3").Value = [double]"0.000549600013252206"
$ws.Range("A4").Value = [double]"7.48244161542317e-05"
$ws.Range("B4").Value = [double]"0.00125667437512966"
$ws.Range("C4").Value = [double]"2.08244792301985e-13"
$ws.Range("D4").Value = [double]"0.0088232506044622"
$ws.Range("E4").Value = [double]"7.50204081632837e-08"
$ws.Range("F4").Value = [double]"0.00118200002665783"
$ws.Range("A5").Value = [double]"4.46602910627902e-05"
$ws.Range("B5").Value = [double]"0.00193227376593692"
$ws.Range("C5").Value = [double]"5.57288973169373e-13"
$ws.Range("D5").Value = [double]"0.00809753297177881"
$ws.Range("E5").Value = [double]"1.93262284290596e-07"
$ws.Range("F5").Value = [double]"0.00188800003974056"
$ws.Range("A6").Value = [double]"3.29426799300425e-05"
$ws.Range("B6").Value = [double]"0.00244128063203818"
$ws.Range("C6").Value = [double]"9.64001882330223e-13"
$ws.Range("D6").Value = [double]"0.00754638834449791"
$ws.Range("E6").Value = [double]"3.31023463932851e-07"
$ws.Range("F6").Value = [double]"0.0024090000477706"
$ws.Range("A7").Value = [double]"2.48799403495712e-05"
$ws.Range("B7").Value = [double]"0.00298280889779496"
$ws.Range("C7").Value = [double]"1.56781797198746e-12"
$ws.Range("D7").Value = [double]"0.00696365558171164"
$ws.Range("E7").Value = [double]"5.35520493396148e-07"
$ws.Range("F7").Value = [double]"0.0029590000552358"
$ws.Range("A8").Value = [double]"1.2550110653497e-05"
$ws.Range("B8").Value = [double]"0.00451333733148972"
$ws.Range("C8").Value = [double]"4.73097696589186e-12"
$ws.Range("D8").Value = [double]"0.00531505628129388"
$ws.Range("E8").Value = [double]"1.60638721640046e-06"
$ws.Range("F8").Value = [double]"0.00450400006952423"
$ws.Range("A9").Value = [double]"9.68725679601052e-06"
$ws.Range("B9").Value = [double]"0.00512396188794407"
$ws.Range("C9").Value = [double]"6.96601368405146e-12"
$ws.Range("D9").Value = [double]"0.00465767543111297"
$ws.Range("E9").Value = [double]"2.36268094296218e-06"
$ws.Range("F9").Value = [double]"0.00511900007246139"
$ws.Range("A10").Value = [double]"7.8939927477496e-06"
$ws.Range("B10").Value = [double]"0.00559755919711815"
$ws.Range("C10").Value = [double]"9.34503620840023e-12"
$ws.Range("D10").Value = [double]"0.00414627340973957"
$ws.Range("E10").Value = [double]"3.16739314228429e-06"
$ws.Range("F10").Value = [double]"0.00559600007335086"
$ws.Range("A11").Value = [double]"6.22906538852753e-06"
$ws.Range("B11").Value = [double]"0.00612344685776672"
$ws.Range("C11").Value = [double]"1.29644503190764e-11"
$ws.Range("D11").Value = [double]"0.00357916204490461"
$ws.Range("E11").Value = [double]"4.39109732867863e-06"
$ws.Range("F11").Value = [double]"0.0061260000728142"

# --- emf_calc_abs_errors: C2:L2 (calc_emf row, shared with rel-errors sheet) and C3:L3 (abs errors) ---
$ws = $wb.Worksheets.Item("emf_calc_abs_errors")
$ws.Range("C2").Value = [double]"162.835364916109"
$ws.Range("D2").Value = [double]"140.562199060342"
$ws.Range("E2").Value = [double]"122.048237822434"
$ws.Range("F2").Value = [double]"107.578935213466"
$ws.Range("G2").Value = [double]"99.0465163073831"
$ws.Range("H2").Value = [double]"91.1760330537543"
$ws.Range("I2").Value = [double]"71.9886861494901"
$ws.Range("J2").Value = [double]"64.7291244515028"
$ws.Range("K2").Value = [double]"58.9894764414753"
$ws.Range("L2").Value = [double]"52.3479417673376"
$ws.Range("C3").Value = [double]"-0.864635083891471"
$ws.Range("D3").Value = [double]"-1.4378009396579"
$ws.Range("E3").Value = [double]"-1.25176217756572"
$ws.Range("F3").Value = [double]"-3.12106478653401"
$ws.Range("G3").Value = [double]"-3.35348369261688"
$ws.Range("H3").Value = [double]"-3.02396694624572"
$ws.Range("I3").Value = [double]"-0.211313850509924"
$ws.Range("J3").Value = [double]"1.22912445150278"
$ws.Range("K3").Value = [double]"2.98947644147531"
$ws.Range("L3").Value = [double]"8.14794176733763"

# --- emf_calc_rel_errors: C2:L2 (calc_emf row, shared) and C3:L3 (rel errors) ---
$ws = $wb.Worksheets.Item("emf_calc_rel_errors")
$ws.Range("C2").Value = [double]"162.835364916109"
$ws.Range("D2").Value = [double]"140.562199060342"
$ws.Range("E2").Value = [double]"122.048237822434"
$ws.Range("F2").Value = [double]"107.578935213466"
$ws.Range("G2").Value = [double]"99.0465163073831"
$ws.Range("H2").Value = [double]"91.1760330537543"
$ws.Range("I2").Value = [double]"71.9886861494901"
$ws.Range("J2").Value = [double]"64.7291244515028"
$ws.Range("K2").Value = [double]"58.9894764414753"
$ws.Range("L2").Value = [double]"52.3479417673376"
$ws.Range("C3").Value = [double]"-0.00528182702438284"
$ws.Range("D3").Value = [double]"-0.0101253587299852"
$ws.Range("E3").Value = [double]"-0.010152166890233"
$ws.Range("F3").Value = [double]"-0.0281939005106956"
$ws.Range("G3").Value = [double]"-0.0327488641857117"
$ws.Range("H3").Value = [double]"-0.0321015599389142"
$ws.Range("I3").Value = [double]"-0.00292678463310145"
$ws.Range("J3").Value = [double]"0.019356290574847"
$ws.Range("K3").Value = [double]"0.0533835078834877"
$ws.Range("L3").Value = [double]"0.184342573921666"

